$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds "K" values. Regenerate them to use K instead of Strike#.
$ws.Range("G2").Value = 6
$ws.Range("G3").Value = 5
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 3
